$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted before what used to be row 7.
# This pushes the old row 7 -> row 8 and the old row 8 -> row 9.
$ws.Rows("7:7").Insert()

# Populate the newly inserted row 7 with the latest "Espárragos" record.
$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "Vega Monumental Concepción"
$ws.Range("C7").Value = "Bíobío"
$ws.Range("D7").Value = 44524
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 300000000
$ws.Range("G7").Value = "Espárragos"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 1600
$ws.Range("M7").Value = 1550
$ws.Range("N7").Value = "$/kilo"
$ws.Range("O7").Value = "Provincia de Talca"
$ws.Range("P7").Value = 1550
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
